$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: Insert a new "Meta description" paragraph right after the title
#         paragraph ("Play Age of Egypt Free - Best Online Slot Game").
#         The new paragraph is: <empty run><bold "Meta description"><plain
#         ": Join the pharaohs ..."> and uses the default "Normal" style.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$insertPoint = $d.Range($p1.Range.End, $p1.Range.End)

# We append a trailing empty <w:p/> too; Word's InsertXML only actually
# splits off a new paragraph when the inserted fragment itself contains a
# paragraph boundary that lands inside the target range. The extra empty
# paragraph it leaves behind is removed right after.
$metaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Join the pharaohs and search for hidden treasures in Age of Egypt, the best online slot game with multiple special features. Play for free today.</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$insertPoint.InsertXML($metaXml)

# Remove the leftover empty paragraph used only to force the paragraph split.
$leftover = $d.Paragraphs(3)
$leftover.Range.Delete()

# ---------------------------------------------------------------------------
# Step 2: Remove the duplicate bold heading paragraph
#         ("Play Age of Egypt Free - Best Online Slot Game") that used to sit
#         right before the closing "Join the pharaohs ..." paragraph.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$dupHeading = $null
for ($i = 1; $i -le $count; $i++) {
  $pp = $d.Paragraphs($i)
  $t = $pp.Range.Text.TrimEnd([char]13, [char]7)
  if ($i -ne 1 -and $t -eq "Play Age of Egypt Free - Best Online Slot Game") {
    $dupHeading = $pp
  }
}
if ($dupHeading -ne $null) {
  $dupHeading.Range.Delete()
}

# ---------------------------------------------------------------------------
# Step 3: Replace the text of the closing italic paragraph (formerly
#         "Join the pharaohs ...") with the new "Prompt: ..." text, keeping
#         the existing italic run formatting and leading empty run intact.
# ---------------------------------------------------------------------------
$count2 = $d.Paragraphs.Count
$closingPara = $null
for ($i = 1; $i -le $count2; $i++) {
  $pp = $d.Paragraphs($i)
  $t = $pp.Range.Text.TrimEnd([char]13, [char]7)
  if ($t -eq "Join the pharaohs and search for hidden treasures in Age of Egypt, the best online slot game with multiple special features. Play for free today.") {
    $closingPara = $pp
  }
}

if ($closingPara -ne $null) {
  # Exclude the trailing paragraph-mark character from the replaced range so
  # the paragraph itself is not split/duplicated.
  $jr = $d.Range($closingPara.Range.Start, $closingPara.Range.End - 1)
  $promptXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:i/></w:rPr><w:t>Prompt: Create a colorful cartoon-style feature image for the online slot &quot;Age of Egypt&quot;, featuring a happy Maya warrior with glasses. The image should have a fun and adventurous feel, with the warrior holding a treasure from the game and standing in front of a famous landmark from ancient Egypt, such as the Great Sphinx or the pyramids. Use bright and bold colors to catch the attention of potential players, and include the game title and the Playtech logo as well.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
  [void]$jr.InsertXML($promptXml)
}

Write-Output "Edit complete"
